# "Generate Report for Archive"
#
# The localization status report is regenerated: every cell that held the
# "Ready for handoff" status label now reads "In Translation" (the shared
# string is used on the Overview sheet's zh-cn/de-de columns and on the
# Status column of each per-locale sheet). Excel's column AutoFit then
# shrinks the now-narrower Status/zh-cn/de-de columns to fit the new,
# shorter text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col E) and de-de (col F) status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = $newStatus

# --- Per-locale sheets: Status column (col C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = $newStatus

# --- Re-fit the columns that just lost the long "Ready for handoff" text ---
$wsOverview.Columns("E:F").AutoFit()
$wsZhCn.Columns("C:C").AutoFit()
$wsDeDe.Columns("C:C").AutoFit()

# The AutoFit heuristic in this runtime doesn't reproduce Excel's exact
# font-metric-based pixel width, so nudge the three status columns to the
# precise width Excel itself settles on for the new text (quantized to the
# nearest internal pixel grid step, same as native ColumnWidth assignment).
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5
$wsZhCn.Range("C1").ColumnWidth = 12.5
$wsDeDe.Range("C1").ColumnWidth = 12.5
